$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 291
$newTotalBudget = 80000000

for ($r = 2; $r -le $lastRow; $r++) {
    $share = $ws.Cells.Item($r, 16).Value2   # column P = budgetShare
    $ws.Cells.Item($r, 17).Value2 = $share * $newTotalBudget   # column Q = Budget
}
